$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9937563538551331
$ws.Range("B1").Value = 1.917658448219299
$ws.Range("C1").Value = 5.432667255401611
$ws.Range("D1").Value = 2.310971021652222
$ws.Range("E1").Value = 1.298062443733215
